# Updated queries for C3DC first half testcases.
#
# The workbook's SQL queries (stored in cells B2, C2, B3, B4, B5, B6, B7 of
# Sheet1) joined df_study/df_participant/etc. using the generic ".id" /
# "study.id" / "participant.id" columns. They need to be updated to use the
# explicit, renamed key columns: "study_id" / "study.study_id" and
# "participant_id" / "participant.participant_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All seven SQL-query cells that contain the JOIN clauses to update.
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $cells) {
    $rng = $ws.Range($cellAddr)
    $query = $rng.Value()

    # Quoted, dotted forms first (more specific), then the bare "x.id" forms.
    $query = $query.Replace('"study.id"', '"study.study_id"')
    $query = $query.Replace('"participant.id"', '"participant.participant_id"')
    $query = $query.Replace('std.id', 'std.study_id')
    $query = $query.Replace('prt.id', 'prt.participant_id')

    $rng.Value = $query
}

# Update the sheet's current selection to match the saved view (C7).
$ws.Range("C7").Select()
